$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Border($rng, $idx, $style, $weight, $color) {
    $b = $rng.Borders.Item($idx)
    $b.LineStyle = $style
    $b.Weight = $weight
    if ($color -ne $null) { $b.Color = $color }
}

function Set-Box($rng, $leftColor, $rightColor, $topColor, $bottomColor) {
    Set-Border $rng 7 1 -4138 $leftColor
    Set-Border $rng 10 1 -4138 $rightColor
    Set-Border $rng 8 1 -4138 $topColor
    Set-Border $rng 9 1 -4138 $bottomColor
}

# ---- Set cell values ----
$ws.Range('A1').Value = 'PLOS articles citation'
$ws.Range('A3').Value = 'Article title'
$ws.Range('B3').Value = 'Authors'
$ws.Range('C3').Value = 'Department'
$ws.Range('D3').Value = 'Citation'
$ws.Range('E3').Value = 'Published'
$ws.Range('A4').Value = 'Truth, Probability, and Frameworks'
$ws.Range('B4').Value = 'Jonathan D Wren'
$ws.Range('D4').Value = 'Wren JD (2005) Truth, Probability, and Frameworks. PLoS Med 2(11): e361. doi:10.1371/journal.pmed.0020361'
$ws.Range('E4').Value = 38685
$ws.Range('A5').Value = 'Plant-Symbiotic Fungi as Chemical Engineers: Multi-Genome Analysis of the Clavicipitaceae Reveals Dynamics of Alkaloid Loci'
$ws.Range('B5').Value = 'Bruce A. Roe'
$ws.Range('C5').Value = 'Department of Chemistry and Biochemistry'
$ws.Range('D5').Value = 'Schardl CL, Young CA, Hesse U, Amyotte SG, Andreeva K, et al. (2013) Plant-Symbiotic Fungi as Chemical Engineers: Multi-Genome Analysis of the Clavicipitaceae Reveals Dynamics of Alkaloid Loci. PLoS Genet 9(2): e1003323. doi:10.1371/journal.pgen.1003323'
$ws.Range('E5').Value = 41333
$ws.Range('A6').Value = 'Cryptocephal, the Drosophila melanogaster ATF4, Is a Specific Coactivator for Ecdysone Receptor Isoform B2'
$ws.Range('B6').Value = 'Sebastien A. Gauthier'
$ws.Range('C6').Value = 'Department of Biolog'
$ws.Range('D6').Value = 'Gauthier SA, VanHaaften E, Cherbas L, Cherbas P, Hewes RS (2012) Cryptocephal, the Drosophila melanogaster ATF4, Is a Specific Coactivator for Ecdysone Receptor Isoform B2. PLoS Genetics 8(8): e1002883. doi:10.1371/journal.pgen.1002883'
$ws.Range('E6').Value = 41130
$ws.Range('B7').Value = 'Randall S. Hewes'
$ws.Range('C7').Value = 'Department of Biolog'
$ws.Range('E7').Value = 41130
$ws.Range('A8').Value = 'Genetic Evidence for an Indispensable Role of Somatic Embryogenesis Receptor Kinases in Brassinosteroid Signaling'
$ws.Range('B8').Value = 'Kai He'
$ws.Range('C8').Value = 'Department of Botany and Microbiology'
$ws.Range('D8').Value = 'Gou X, Yin H, He K, Du J, Yi J, et al. (2012) Genetic Evidence for an Indispensable Role of Somatic Embryogenesis Receptor Kinases in Brassinosteroid Signaling. PLoS Genet 8(1): e1002452. doi:10.1371/journal.pgen.1002452'
$ws.Range('E8').Value = 40920
$ws.Range('B9').Value = 'Jia Li'
$ws.Range('C9').Value = 'Department of Botany and Microbiology'
$ws.Range('E9').Value = 40920
$ws.Range('A10').Value = 'The Thermoanaerobacter Glycobiome Reveals Mechanisms of Pentose and Hexose Co-Utilization in Bacteria'
$ws.Range('B10').Value = 'Lu Lin'
$ws.Range('C10').Value = 'Department of Botany and Microbiology'
$ws.Range('D10').Value = 'Lin L, Song H, Tu Q, Qin Y, Zhou A, et al. (2011) The Thermoanaerobacter Glycobiome Reveals Mechanisms of Pentose and Hexose Co-Utilization in Bacteria. PLoS Genet 7(10): e1002318. doi:10.1371/journal.pgen.1002318'
$ws.Range('E10').Value = 40829
$ws.Range('B11').Value = 'Qichao Tu'
$ws.Range('C11').Value = 'Department of Botany and Microbiology'
$ws.Range('E11').Value = 40829
$ws.Range('B12').Value = 'Yujia Qin'
$ws.Range('C12').Value = 'Department of Botany and Microbiology'
$ws.Range('E12').Value = 40829
$ws.Range('B13').Value = 'Aifen Zhou'
$ws.Range('C13').Value = 'Department of Botany and Microbiology'
$ws.Range('E13').Value = 40829
$ws.Range('B14').Value = 'Wenbin Liu'
$ws.Range('C14').Value = 'Department of Botany and Microbiology'
$ws.Range('E14').Value = 40829
$ws.Range('B15').Value = 'Zhili He'
$ws.Range('C15').Value = 'Department of Botany and Microbiology'
$ws.Range('E15').Value = 40829
$ws.Range('B16').Value = 'Jizhong Zhou'
$ws.Range('C16').Value = 'Department of Botany and Microbiology'
$ws.Range('E16').Value = 40829
$ws.Range('A17').Value = 'A Genome-Wide Survey of Imprinted Genes in Rice Seeds Reveals Imprinting Primarily Occurs in the Endosperm'
$ws.Range('B17').Value = 'Scott Russell'
$ws.Range('C17').Value = 'Department of Botany and Microbiology'
$ws.Range('D17').Value = 'Luo M, Taylor JM, Spriggs A, Zhang H, Wu X, et al. (2011) A Genome-Wide Survey of Imprinted Genes in Rice Seeds Reveals Imprinting Primarily Occurs in the Endosperm. PLoS Genet 7(6): e1002125. doi:10.1371/journal.pgen.1002125'
$ws.Range('E17').Value = 40717
$ws.Range('A18').Value = 'Spatial Distribution and Risk Factors of Highly Pathogenic Avian Influenza (HPAI) H5N1 in China'
$ws.Range('B18').Value = 'Xiangming Xiao'
$ws.Range('C18').Value = 'Department of Botany and Microbiology'
$ws.Range('D18').Value = 'Martin V, Pfeiffer DU, Zhou X, Xiao X, Prosser DJ, et al. (2011) Spatial Distribution and Risk Factors of Highly Pathogenic Avian Influenza (HPAI) H5N1 in China. PLoS Pathog 7(3): e1001308. doi:10.1371/journal.ppat.1001308'
$ws.Range('E18').Value = 40605
$ws.Range('A19').Value = 'Israeli Acute Paralysis Virus: Epidemiology, Pathogenesis and Implications for Honey Bee Health'
$ws.Range('B19').Value = 'Ai Fen Zhou'
$ws.Range('C19').Value = 'Institute for Environmental Genomics (IEG), University of Oklahoma'
$ws.Range('D19').Value = 'Chen YP, Pettis JS, Corona M, Chen WP, Li CJ, et al. (2014) Israeli Acute Paralysis Virus: Epidemiology, Pathogenesis and Implications for Honey Bee Health. PLoS Pathog 10(7): e1004261. doi:10.1371/journal.ppat.1004261'
$ws.Range('E19').Value = 41851
$ws.Range('B20').Value = 'Li You Wu'
$ws.Range('B21').Value = 'Ji Zhong Zhou'
$ws.Range('A22').Value = 'Sorghum Genome Sequencing by Methylation Filtration'
$ws.Range('B22').Value = 'Graham Wiley'
$ws.Range('C22').Value = 'Department of Chemistry and Biochemistry'
$ws.Range('E22').Value = 38356
$ws.Range('B23').Value = 'Bruce A Roe'
$ws.Range('C23').Value = 'Department of Chemistry and Biochemistry'
$ws.Range('D23').Value = 'Bedell JA, Budiman MA, Nunberg A, Citek RW, Robbins D, et al. (2005) Sorghum Genome Sequencing by Methylation Filtration. PLoS Biol 3(1): e13. doi:10.1371/journal.pbio.0030013'
$ws.Range('E23').Value = 38356

# ---- Clear cells that should be empty in target but had values before ----
$ws.Range('C4').ClearContents()
$ws.Range('F3').Style = 'Normal'
